# Shift the "August" schedule table (rows 6-36, columns A/B) 22 days earlier.
# Column A holds the German weekday label, column B the serial date. Both are
# recomputed together so the weekday label always matches the new date's
# position in the existing 7-day (Mittwoch/Donnerstag/Freitag/Samstag/
# Sonntag/Montag/Dienstag) rotation used throughout this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$weekdays = @("Mittwoch", "Donnerstag", "Freitag", "Samstag", "Sonntag", "Montag", "Dienstag")

$firstRow = 6
$lastRow = 36
$firstDate = 42194

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $offset = $row - $firstRow
    $dayName = $weekdays[$offset % 7]
    $dateSerial = $firstDate + $offset

    $ws.Range("A$row").Value = $dayName
    $ws.Range("B$row").Value = $dateSerial
}
